$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric stat corrections ---
$ws.Range("BC2").Value = 18
$ws.Range("D4").Value = 80
$ws.Range("F4").Value = 36
$ws.Range("G4").Value = 0.55
$ws.Range("I4").Value = 35.9
$ws.Range("J4").Value = 78
$ws.Range("L4").Value = 8.6
$ws.Range("N4").Value = 0.369
$ws.Range("Q4").Value = 0.754
$ws.Range("S4").Value = 29.4
$ws.Range("V4").Value = 14.4
$ws.Range("X4").Value = 3.9
$ws.Range("AA4").Value = 20.7
$ws.Range("AC4").Value = -0.5
$ws.Range("AD4").Value = 27
$ws.Range("AP4").Value = 11
$ws.Range("AV4").Value = 11
$ws.Range("BC4").Value = 17
$ws.Range("AH5").Value = 7
$ws.Range("AP5").Value = 10
$ws.Range("AI7").Value = 20
$ws.Range("D9").Value = 80
$ws.Range("F9").Value = 44
$ws.Range("G9").Value = 0.45
$ws.Range("J9").Value = 86
$ws.Range("M9").Value = 23.7
$ws.Range("N9").Value = 0.361
$ws.Range("O9").Value = 18.8
$ws.Range("P9").Value = 26.1
$ws.Range("Q9").Value = 0.721
$ws.Range("S9").Value = 33.2
$ws.Range("V9").Value = 15.8
$ws.Range("AC9").Value = -2
$ws.Range("AD9").Value = 27
$ws.Range("AF9").Value = 18
$ws.Range("AH9").Value = 25
$ws.Range("AI9").Value = 11
$ws.Range("AZ11").Value = 22
$ws.Range("AW12").Value = 16
$ws.Range("D14").Value = 80
$ws.Range("E14").Value = 56
$ws.Range("G14").Value = 0.7
$ws.Range("I14").Value = 39.1
$ws.Range("N14").Value = 0.354
$ws.Range("O14").Value = 21.2
$ws.Range("AB14").Value = 107.9
$ws.Range("AD14").Value = 27
$ws.Range("AV14").Value = 7
$ws.Range("AP16").Value = 28
$ws.Range("AH17").Value = 7
$ws.Range("AI17").Value = 12
$ws.Range("AH18").Value = 7
$ws.Range("AJ18").Value = 20
$ws.Range("AV19").Value = 8
$ws.Range("AH20").Value = 11
$ws.Range("D21").Value = 80
$ws.Range("E21").Value = 35
$ws.Range("G21").Value = 0.438
$ws.Range("I21").Value = 37
$ws.Range("J21").Value = 82.40000000000001
$ws.Range("N21").Value = 0.37
$ws.Range("O21").Value = 15.3
$ws.Range("P21").Value = 20.1
$ws.Range("Q21").Value = 0.76
$ws.Range("U21").Value = 20.2
$ws.Range("AB21").Value = 98.5
$ws.Range("AC21").Value = -1
$ws.Range("AD21").Value = 27
$ws.Range("AE21").Value = 20
$ws.Range("AF21").Value = 20
$ws.Range("AG21").Value = 20
$ws.Range("AJ21").Value = 19
$ws.Range("AO21").Value = 29
$ws.Range("AP21").Value = 30
$ws.Range("AQ21").Value = 14
$ws.Range("AW21").Value = 17
$ws.Range("AI23").Value = 20
$ws.Range("AH25").Value = 26
$ws.Range("AH26").Value = 11
$ws.Range("AQ27").Value = 13
$ws.Range("AH28").Value = 26
$ws.Range("AP28").Value = 29
$ws.Range("AV28").Value = 12
$ws.Range("AZ29").Value = 30
$ws.Range("AO31").Value = 28

# --- Date column fix: "4-15-2013-14" -> "2014-04-15" (text, not a date) ---
$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).Value = "2014-04-15"
}
$dateRange.Style = "Normal"
